$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Korean verb vocabulary rows 32-101 (extends the existing A1:B31 table)
$data = @(
    @("열다", "to open"),
    @("닫다", "to close"),
    @("입다", "to wear"),
    @("벗다", "to take off"),
    @("신다", "to wear (shoes/socks)"),
    @("타다", "to ride"),
    @("내리다", "to get off"),
    @("걷다", "to walk"),
    @("뛰다", "to run"),
    @("앉다", "to sit"),
    @("서다", "to stand"),
    @("알다", "to know"),
    @("모르다", "to not know"),
    @("생각하다", "to think"),
    @("기억하다", "to remember"),
    @("잊다", "to forget"),
    @("느끼다", "to feel"),
    @("보내다", "to send"),
    @("빌리다", "to borrow"),
    @("돌려주다", "to return"),
    @("돕다", "to help"),
    @("필요하다", "to need"),
    @("찾다", "to find/look for"),
    @("잃어버리다", "to lose"),
    @("시작하다", "to start"),
    @("끝나다", "to end"),
    @("기다리다", "to wait"),
    @("준비하다", "to prepare"),
    @("약속하다", "to promise"),
    @("거짓말하다", "to lie"),
    @("고백하다", "to confess"),
    @("결혼하다", "to marry"),
    @("이혼하다", "to divorce"),
    @("태어나다", "to be born"),
    @("자라다", "to grow up"),
    @("이기다", "to win"),
    @("지다", "to lose (game)"),
    @("싸우다", "to fight"),
    @("화내다", "to get angry"),
    @("용서하다", "to forgive"),
    @("사과하다", "to apologize"),
    @("감사하다", "to thank"),
    @("축하하다", "to congratulate"),
    @("노래하다", "to sing"),
    @("춤추다", "to dance"),
    @("그리다", "to draw"),
    @("만들다", "to make"),
    @("고치다", "to fix"),
    @("청소하다", "to clean"),
    @("요리하다", "to cook"),
    @("빨래하다", "to do laundry"),
    @("운전하다", "to drive"),
    @("여행하다", "to travel"),
    @("운동하다", "to exercise"),
    @("공부하다", "to study"),
    @("일하다", "to work"),
    @("쉬다", "to rest"),
    @("잠자다", "to sleep"),
    @("꿈꾸다", "to dream"),
    @("믿다", "to believe"),
    @("바라다", "to hope/wish"),
    @("걱정하다", "to worry"),
    @("기뻐하다", "to rejoice"),
    @("슬퍼하다", "to be sad"),
    @("놀라다", "to be surprised"),
    @("무서워하다", "to be afraid"),
    @("부러워하다", "to envy"),
    @("질투하다", "to be jealous"),
    @("존경하다", "to respect"),
    @("무시하다", "to ignore"),
)

$startRow = 32
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

